# This script applies the change described by the diff to sheet1:
#  - Row 3 ("CROANCA"/"Facebook" garbage row) is removed entirely, shifting
#    all subsequent rows up by one (dimension goes from A1:G12 to A1:G11).
#  - A handful of numeric values in the resulting rows 3, 4, 5 and 10 are
#    corrected to their new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old row 3 (CROANCA/Facebook) - this shifts rows 4-12 up to 3-11.
$ws.Rows.Item(3).Delete()

# Now apply the corrected values for the rows that changed beyond the plain shift.
# New row 3: CRONACA / Facebook
$ws.Cells.Item(3, 4).Value = 15
$ws.Cells.Item(3, 5).Value = 1093
$ws.Cells.Item(3, 6).Value = 85
$ws.Cells.Item(3, 7).Value = 3

# New row 4: CRONACA / Instagram
$ws.Cells.Item(4, 5).Value = 995
$ws.Cells.Item(4, 6).Value = 181

# New row 5: CRONACA / YouTube
$ws.Cells.Item(5, 5).Value = 1101

# New row 10: POLITICA / Instagram
$ws.Cells.Item(10, 5).Value = 1008

Write-Output "Done"
